$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-02-21 Wednesday" "2024-02-22 Thursday"

Replace-Text "324×4=" "911×2="
Replace-Text "748×5=" "559×8="
Replace-Text "679×6=" "861×8="
Replace-Text "572×4=" "812×9="
Replace-Text "707×5=" "622×7="
Replace-Text "239×6=" "942×5="
Replace-Text "370×3=" "808×2="
Replace-Text "529×3=" "868×8="
Replace-Text "376×6=" "925×4="
Replace-Text "443×4=" "500×2="
Replace-Text "194×9=" "211×5="
Replace-Text "610×5=" "828×8="
Replace-Text "633×3=" "216×4="
Replace-Text "670×4=" "736×2="
Replace-Text "772×4=" "949×3="
Replace-Text "394×2=" "743×8="
Replace-Text "431×7=" "880×8="
Replace-Text "877×2=" "212×3="
Replace-Text "594×2=" "662×7="
Replace-Text "676×6=" "928×5="
Replace-Text "952×9=" "190×6="
Replace-Text "892×8=" "558×9="
Replace-Text "697×7=" "863×4="
Replace-Text "495×8=" "510×2="
Replace-Text "822×4=" "831×4="
